$d = $word.ActiveDocument

# Locate the paragraph that currently holds "int numeroVidas;" (split across
# runs "i" / "nt" / " " / "numeroVidas" / ";" with the _GoBack bookmark sitting
# between "i" and "nt"). We identify it by its visible text rather than a
# hard-coded index, so the script is resilient to minor structural changes.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "int numeroVidas;*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph containing 'int numeroVidas;'"
}

$p = $d.Paragraphs($targetIndex)

# Replace the whole paragraph (including the embedded _GoBack bookmark) with
# two new paragraphs:
#   1) "String color;" (new line for the cat's color field), keeping the
#      _GoBack bookmark at the end, right where it used to sit.
#   2) "int numeroVidas;" rebuilt with clean spell-check proofErr markers,
#      same as the other declarations in the document.
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>String</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> color;</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>int</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>numeroVidas</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>;</w:t></w:r>
</w:p>
"@

$p.Range.InsertXML($xml)
